$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 261.3203778131603
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 1396.132104504975
